$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for all data rows (2-54)
# from 45224 (2023-10-25) to 45233 (2023-11-03)
$ws.Range("C2:C54").Value = 45233
